$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 199.55556
$ws.Range("I42").Value = 155
$ws.Range("J42").Value = 221.83333
$ws.Range("K42").Value = 465
$ws.Range("L42").Value = 665.49999
$ws.Range("M42").Value = -235
$ws.Range("N42").Value = -1125.49999

$ws.Range("H111").Value = 8437.5
$ws.Range("I111").Value = 8437.5
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 25312.5
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -22245.5
$ws.Range("N111").ClearContents()

$ws.Range("H132").Value = 24751.895
$ws.Range("I132").Value = 3810.524
$ws.Range("J132").Value = 200659.4
$ws.Range("K132").Value = 11431.572
$ws.Range("L132").Value = 601978.2
$ws.Range("M132").Value = -8901.572
$ws.Range("N132").Value = -607038.2

$ws.Range("H138").Value = 5274.3726
$ws.Range("I138").Value = 6061.923
$ws.Range("J138").Value = 5004.9473
$ws.Range("K138").Value = 18185.769
$ws.Range("L138").Value = 15014.8419
$ws.Range("M138").Value = -13045.769
$ws.Range("N138").Value = -25294.8419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1324255.4
$ws.Range("I32").Value = 589417.4
$ws.Range("J32").Value = 47619050
$ws.Range("K32").Value = 589417.4
$ws.Range("L32").Value = 47619050
$ws.Range("M32").Value = -589130.4

$ws.Range("H45").Value = 62553784
$ws.Range("I45").Value = 70876.5
$ws.Range("J45").Value = 250002510
$ws.Range("K45").Value = 70876.5
$ws.Range("L45").Value = 250002510
$ws.Range("M45").Value = -70499.5
$ws.Range("N45").Value = -250003264

$ws.Range("H122").Value = 3854.5625
$ws.Range("I122").Value = 3074.889
$ws.Range("J122").Value = 4857
$ws.Range("K122").Value = 9224.667000000001
$ws.Range("L122").Value = 14571
$ws.Range("M122").Value = -6774.667000000001

$ws.Range("H132").Value = 3265.9473
$ws.Range("I132").Value = 2623.2173
$ws.Range("J132").Value = 4251.467
$ws.Range("K132").Value = 7869.651899999999
$ws.Range("L132").Value = 12754.401
$ws.Range("M132").Value = -5339.651899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2477.8462
$ws.Range("I134").Value = 2396.2
$ws.Range("J134").Value = 2750
$ws.Range("K134").Value = 7188.599999999999
$ws.Range("L134").Value = 8250
$ws.Range("M134").Value = -4653.599999999999
$ws.Range("N134").Value = -13320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2362438.5
$ws.Range("I31").Value = 2700.7144
$ws.Range("J31").Value = 3209523.8
$ws.Range("K31").Value = 2700.7144
$ws.Range("L31").Value = 3209523.8
$ws.Range("M31").Value = -2405.7144
$ws.Range("N31").Value = -3210113.8

$ws.Range("H34").Value = 2362438.5
$ws.Range("I34").Value = 2700.7144
$ws.Range("J34").Value = 3209523.8
$ws.Range("K34").Value = 2700.7144
$ws.Range("L34").Value = 3209523.8
$ws.Range("M34").Value = -2498.7144
$ws.Range("N34").Value = -3209927.8

$ws.Range("H132").Value = 2718.5652
$ws.Range("I132").Value = 2718.5652
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8155.6956
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5625.6956
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1644538.2
$ws.Range("I68").Value = 4626.6924
$ws.Range("J68").Value = 2088681
$ws.Range("K68").Value = 13880.0772
$ws.Range("L68").Value = 6266043
$ws.Range("M68").Value = -13069.0772
$ws.Range("N68").Value = -6267665

$ws.Range("H71").Value = 1644538.2
$ws.Range("I71").Value = 4626.6924
$ws.Range("J71").Value = 2088681
$ws.Range("K71").Value = 41640.2316
$ws.Range("L71").Value = 18798129
$ws.Range("M71").Value = -37584.2316
$ws.Range("N71").Value = -18806241

$ws.Range("H113").Value = 809.9167
$ws.Range("I113").Value = 831.25
$ws.Range("J113").Value = 805.65
$ws.Range("K113").Value = 2493.75
$ws.Range("L113").Value = 2416.95
$ws.Range("M113").Value = -323.75
$ws.Range("N113").Value = -6756.95

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 142859740
$ws.Range("I80").Value = 250001420
$ws.Range("J80").Value = 4166.3335
$ws.Range("K80").Value = 250001420
$ws.Range("L80").Value = 4166.3335
$ws.Range("M80").Value = -250000422
$ws.Range("N80").Value = -6162.3335

$ws.Range("H83").Value = 142859740
$ws.Range("I83").Value = 250001420
$ws.Range("J83").Value = 4166.3335
$ws.Range("K83").Value = 1250007100
$ws.Range("L83").Value = 20831.6675
$ws.Range("M83").Value = -1250002108
$ws.Range("N83").Value = -30815.6675

$ws.Range("H122").Value = 28574838
$ws.Range("I122").Value = 3012.72
$ws.Range("J122").Value = 100004400
$ws.Range("K122").Value = 9038.16
$ws.Range("L122").Value = 300013200
$ws.Range("M122").Value = -6588.16
$ws.Range("N122").Value = -300018100

$ws.Range("H132").Value = 1335.238
$ws.Range("I132").Value = 1279.3889
$ws.Range("J132").Value = 1670.3334
$ws.Range("K132").Value = 3838.1667
$ws.Range("L132").Value = 5011.0002
$ws.Range("M132").Value = -1308.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1962.3182
$ws.Range("I7").Value = 1603.5
$ws.Range("J7").Value = 2590.25
$ws.Range("K7").Value = 1603.5
$ws.Range("L7").Value = 2590.25
$ws.Range("M7").Value = -1491.5

$ws.Range("H40").Value = 63109.477
$ws.Range("I40").Value = 91178.57000000001
$ws.Range("J40").Value = 6971.2856
$ws.Range("K40").Value = 91178.57000000001
$ws.Range("L40").Value = 6971.2856
$ws.Range("M40").Value = -91042.57000000001

$ws.Range("H82").Value = 1572.7273
$ws.Range("I82").Value = 1525
$ws.Range("J82").Value = 1630
$ws.Range("K82").Value = 1525
$ws.Range("L82").Value = 1630
$ws.Range("M82").Value = -1164
$ws.Range("N82").Value = -2352

$ws.Range("H85").Value = 1572.7273
$ws.Range("I85").Value = 1525
$ws.Range("J85").Value = 1630
$ws.Range("K85").Value = 1525
$ws.Range("L85").Value = 1630
$ws.Range("M85").Value = -277
$ws.Range("N85").Value = -4126

$ws.Range("H100").Value = 3862719
$ws.Range("I100").Value = 4505839
$ws.Range("J100").Value = 3999
$ws.Range("K100").Value = 4505839
$ws.Range("L100").Value = 3999
$ws.Range("M100").Value = -4505298
$ws.Range("N100").Value = -5081

$ws.Range("H122").Value = 2845.5454
$ws.Range("I122").Value = 3122.4285
$ws.Range("J122").Value = 2361
$ws.Range("K122").Value = 9367.2855
$ws.Range("L122").Value = 7083
$ws.Range("M122").Value = -6917.2855
$ws.Range("N122").Value = -11983

$ws.Range("H126").Value = 1962.3182
$ws.Range("I126").Value = 1603.5
$ws.Range("J126").Value = 2590.25
$ws.Range("K126").Value = 4810.5
$ws.Range("L126").Value = 7770.75
$ws.Range("M126").Value = -2340.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H30").Value = 20006382
$ws.Range("I30").Value = 5636.3335
$ws.Range("J30").Value = 50007500
$ws.Range("K30").Value = 5636.3335
$ws.Range("L30").Value = 50007500
$ws.Range("M30").Value = -5529.3335
$ws.Range("N30").Value = -50007714

$ws.Range("H122").Value = 50002384
$ws.Range("I122").Value = 904
$ws.Range("J122").Value = 62502750
$ws.Range("K122").Value = 2712
$ws.Range("L122").Value = 187508250
$ws.Range("M122").Value = -262
$ws.Range("N122").Value = -187513150

$ws.Range("H123").Value = 59142.668
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 59142.668
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 59142.668
$ws.Range("N123").Value = -68942.66800000001

$ws.Range("H126").Value = 9633.691999999999
$ws.Range("I126").Value = 10019.833
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 30059.499
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -27589.499

$ws.Range("H132").Value = 1715.5385
$ws.Range("I132").Value = 1788.0605
$ws.Range("J132").Value = 1316.6666
$ws.Range("K132").Value = 5364.181500000001
$ws.Range("L132").Value = 3949.9998
$ws.Range("M132").Value = -2834.181500000001
$ws.Range("N132").Value = -9009.9998
